$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.768.75'
$ws.Range('D3').Value = '1.634.17'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.68'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0633'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.46%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.59'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.50%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0792'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.54%  '
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '1.860.50'
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').Value = '1.633.84'
$ws.Range('E14').Value = '  -0.18%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.559'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.20'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '25.792.43'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('E20').Value = '  +1.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '192.40'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.93%  '
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.28'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.26%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.81'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.85%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '142.54'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.21%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.123'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.71%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.89'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0492'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.33'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('E33').Value = '  -0.95%  '
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').Value = '1.131.32'
$ws.Range('E37').Value = '  +1.79%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.52'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.545'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.21%  '
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.56'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '100.78'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.798'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('D46').Value = '1.769.75'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '55.40'
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0505'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.417'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.56%  '
$ws.Range('E51').Value = '  +3.22%  '

Write-Host "Applied all crypto list updates"
